# The commit swaps the two embedded themes in the deck:
#   ppt/theme/theme1.xml  (the main "Integral" / "Red Violet" design, used by the
#                           slide master) ends up holding the colours that used
#                           to live in ppt/theme/theme2.xml (the default
#                           "Office Theme" palette used by the notes master).
#   ppt/theme/theme2.xml  ends up holding the "Integral" / "Red Violet" colours
#                           that used to live in ppt/theme/theme1.xml.
#
# The font scheme and format scheme are byte-for-byte identical between the two
# themes already, so the only observable difference is the 12-slot colour
# scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink). We reassign every slot on
# the presentation's theme (reachable via the slide master's Design) to the
# colours the *other* theme used to hold.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Colours theme1.xml ("Integral" / "Red Violet") had before the edit -
# these are the values that theme2.xml must end up with.
$oldTheme1 = @(
    (0 + 0*256 + 0*65536),             # 1  dk1      000000
    (255 + 255*256 + 255*65536),       # 2  lt1      FFFFFF
    (0x45 + 0x45*256 + 0x51*65536),    # 3  dk2      454551
    (0xD8 + 0xD9*256 + 0xDC*65536),    # 4  lt2      D8D9DC
    (0xE3 + 0x2D*256 + 0x91*65536),    # 5  accent1  E32D91
    (0xC8 + 0x30*256 + 0xCC*65536),    # 6  accent2  C830CC
    (0x4E + 0xA6*256 + 0xDC*65536),    # 7  accent3  4EA6DC
    (0x47 + 0x75*256 + 0xE7*65536),    # 8  accent4  4775E7
    (0x89 + 0x71*256 + 0xE1*65536),    # 9  accent5  8971E1
    (0xD5 + 0x47*256 + 0x73*65536),    # 10 accent6  D54773
    (0x6B + 0x9F*256 + 0x25*65536),    # 11 hlink    6B9F25
    (0x8C + 0x8C*256 + 0x8C*65536)     # 12 folHlink 8C8C8C
)

# Colours theme2.xml ("Office Theme") had before the edit - these are the
# values that theme1.xml must end up with.
$oldTheme2 = @(
    (0 + 0*256 + 0*65536),             # 1  dk1      000000
    (255 + 255*256 + 255*65536),       # 2  lt1      FFFFFF
    (0x44 + 0x54*256 + 0x6A*65536),    # 3  dk2      44546A
    (0xE7 + 0xE6*256 + 0xE6*65536),    # 4  lt2      E7E6E6
    (0x5B + 0x9B*256 + 0xD5*65536),    # 5  accent1  5B9BD5
    (0xED + 0x7D*256 + 0x31*65536),    # 6  accent2  ED7D31
    (0xA5 + 0xA5*256 + 0xA5*65536),    # 7  accent3  A5A5A5
    (0xFF + 0xC0*256 + 0x00*65536),    # 8  accent4  FFC000
    (0x44 + 0x72*256 + 0xC4*65536),    # 9  accent5  4472C4
    (0x70 + 0xAD*256 + 0x47*65536),    # 10 accent6  70AD47
    (0x05 + 0x63*256 + 0xC1*65536),    # 11 hlink    0563C1
    (0x95 + 0x4F*256 + 0x72*65536)     # 12 folHlink 954F72
)

# theme1.xml (this deck's active/presentation theme) takes on the colours
# that used to be theme2's ("Office Theme").
for ($i = 1; $i -le 12; $i++) {
    $colors.Item($i).RGB = $oldTheme2[$i - 1]
}
